# Update column F ("dSF") values on Sheet1 to reflect repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 3
    3  = -1
    6  = 2
    8  = -7
    11 = 4
    12 = -4
    13 = 3
    14 = 0
    15 = -5
    16 = -1
    17 = -1
    18 = 2
    19 = -3
    20 = 1
    21 = -2
    22 = 1
    23 = -5
    24 = -1
    25 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
